{"js": "// Replace the 100 arithmetic-expression cells in the single table, in row-major\n// order, with their updated values (per the commit's regenerated answer key).\n// Cell (r, c) for r in [0,20), c in [0,5) maps 1:1 to newValues[r*5 + c].\nconst newValues = [\"64+32=96\", \"83-14=69\", \"1+3=4\", \"97-42=55\", \"24+47=71\", \"77-76=1\", \"33-2=31\", \"8+20=28\", \"77+22=99\", \"40+58=98\", \"55-20=35\", \"11+0=11\", \"24+65=89\", \"90-12=78\", \"35-28=7\", \"46-2=44\", \"19+52=71\", \"65+27=92\", \"74+5=79\", \"60-51=9\", \"44+29=73\", \"24-11=13\", \"27+19=46\", \"50-3=47\", \"67-30=37\", \"75-66=9\", \"19-5=14\", \"57-20=37\", \"32+36=68\", \"45-19=26\", \"70-3=67\", \"37+58=95\", \"73-24=49\", \"51+2=53\", \"8+68=76\", \"60+36=96\", \"44-34=10\", \"88-6=82\", \"42-28=14\", \"74-48=26\", \"33+17=50\", \"45-14=31\", \"76-34=42\", \"40+37=77\", \"8+16=24\", \"49+36=85\", \"67-19=48\", \"87-52=35\", \"67-19=48\", \"67+20=87\", \"88-74=14\", \"61-27=34\", \"28-15=13\", \"75-5=70\", \"81-71=10\", \"73-26=47\", \"76-25=51\", \"97-9=88\", \"8+38=46\", \"38-28=10\", \"43+35=78\", \"98-42=56\", \"95-85=10\", \"74-66=8\", \"81-67=14\", \"10+56=66\", \"22+71=93\", \"82-46=36\", \"60+38=98\", \"25-23=2\", \"25-22=3\", \"13+81=94\", \"26+34=60\", \"46+43=89\", \"68-1=67\", \"69+12=81\", \"72-48=24\", \"59-35=24\", \"73-27=46\", \"69+18=87\", \"35-11=24\", \"39+0=39\", \"48-14=34\", \"70-56=14\", \"57-35=22\", \"80-26=54\", \"28-11=17\", \"81-80=1\", \"49+4=53\", \"30+19=49\", \"36-23=13\", \"14-5=9\", \"17+14=31\", \"8+32=40\", \"95-48=47\", \"40+2=42\", \"20-3=17\", \"4+91=95\", \"25+4=29\", \"19-15=4\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 5;\n\nif (rowCount * colCount !== newValues.length) {\n  throw new Error(\n    `Table shape mismatch: expected ${newValues.length} cells, found ${rowCount}x${colCount}.`\n  );\n}\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression cells in the single table, in\n# row-major order, with their updated values (per the commit's regenerated\n# answer key). Cell (r, c) (1-based, Word COM convention) for r in [1,20],\n# c in [1,5] maps 1:1 to $newValues[(r-1)*5 + (c-1)].\n$newValues = @(\n    \"64+32=96\",\n    \"83-14=69\",\n    \"1+3=4\",\n    \"97-42=55\",\n    \"24+47=71\",\n    \"77-76=1\",\n    \"33-2=31\",\n    \"8+20=28\",\n    \"77+22=99\",\n    \"40+58=98\",\n    \"55-20=35\",\n    \"11+0=11\",\n    \"24+65=89\",\n    \"90-12=78\",\n    \"35-28=7\",\n    \"46-2=44\",\n    \"19+52=71\",\n    \"65+27=92\",\n    \"74+5=79\",\n    \"60-51=9\",\n    \"44+29=73\",\n    \"24-11=13\",\n    \"27+19=46\",\n    \"50-3=47\",\n    \"67-30=37\",\n    \"75-66=9\",\n    \"19-5=14\",\n    \"57-20=37\",\n    \"32+36=68\",\n    \"45-19=26\",\n    \"70-3=67\",\n    \"37+58=95\",\n    \"73-24=49\",\n    \"51+2=53\",\n    \"8+68=76\",\n    \"60+36=96\",\n    \"44-34=10\",\n    \"88-6=82\",\n    \"42-28=14\",\n    \"74-48=26\",\n    \"33+17=50\",\n    \"45-14=31\",\n    \"76-34=42\",\n    \"40+37=77\",\n    \"8+16=24\",\n    \"49+36=85\",\n    \"67-19=48\",\n    \"87-52=35\",\n    \"67-19=48\",\n    \"67+20=87\",\n    \"88-74=14\",\n    \"61-27=34\",\n    \"28-15=13\",\n    \"75-5=70\",\n    \"81-71=10\",\n    \"73-26=47\",\n    \"76-25=51\",\n    \"97-9=88\",\n    \"8+38=46\",\n    \"38-28=10\",\n    \"43+35=78\",\n    \"98-42=56\",\n    \"95-85=10\",\n    \"74-66=8\",\n    \"81-67=14\",\n    \"10+56=66\",\n    \"22+71=93\",\n    \"82-46=36\",\n    \"60+38=98\",\n    \"25-23=2\",\n    \"25-22=3\",\n    \"13+81=94\",\n    \"26+34=60\",\n    \"46+43=89\",\n    \"68-1=67\",\n    \"69+12=81\",\n    \"72-48=24\",\n    \"59-35=24\",\n    \"73-27=46\",\n    \"69+18=87\",\n    \"35-11=24\",\n    \"39+0=39\",\n    \"48-14=34\",\n    \"70-56=14\",\n    \"57-35=22\",\n    \"80-26=54\",\n    \"28-11=17\",\n    \"81-80=1\",\n    \"49+4=53\",\n    \"30+19=49\",\n    \"36-23=13\",\n    \"14-5=9\",\n    \"17+14=31\",\n    \"8+32=40\",\n    \"95-48=47\",\n    \"40+2=42\",\n    \"20-3=17\",\n    \"4+91=95\",\n    \"25+4=29\",\n    \"19-15=4\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nif ($rows * $cols -ne $newValues.Count) {\n    throw \"Table shape mismatch: expected $($newValues.Count) cells, found $($rows)x$($cols).\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
